# Refresh the Price (D) and Volume(1h) (E) columns in the cryptos list
# with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.388.96"
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = "'1.848.30"
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('D4').Value = "'0.9985"
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'240.44"
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = "'1.000"
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = "'0.07500"
$ws.Range('D9').Value = "'0.2903"
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').Value = "'24.45"
$ws.Range('E10').Value = '  -1.16%  '
$ws.Range('D11').Value = "'0.07738"
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('D12').Value = "'1.848.44"
$ws.Range('E12').Value = '  -2.11%  '
$ws.Range('D13').Value = "'4.998"
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').Value = "'0.6808"
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').Value = "'0.00001043"
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').Value = "'82.18"
$ws.Range('E16').Value = '  -1.27%  '
$ws.Range('D17').Value = "'2.106.12"
$ws.Range('E17').Value = '  -3.70%  '
$ws.Range('D18').Value = "'6.185"
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').Value = "'29.430.11"
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = "'229.35"
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = "'1.000"
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = "'7.468"
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').Value = "'0.9997"
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').Value = "'0.1377"
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D27').Value = "'8.412"
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = "'17.53"
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('D29').Value = "'0.06442"
$ws.Range('E29').Value = '  +14.90%  '
$ws.Range('D30').Value = "'1.390"
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').Value = "'1.475"
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').Value = "'4.095"
$ws.Range('E32').Value = '  -0.55%  '
$ws.Range('D33').Value = "'4.068"
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').Value = "'1.829"
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('D35').Value = "'1.140"
$ws.Range('E35').Value = '  -1.86%  '
$ws.Range('D36').Value = "'0.6989"
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('D37').Value = "'2.582"
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').Value = "'1.264.35"
$ws.Range('E38').Value = '  +2.79%  '
$ws.Range('D39').Value = "'2.837"
$ws.Range('E39').Value = '  +4.52%  '
$ws.Range('D40').Value = "'0.01830"
$ws.Range('E40').Value = '  +1.57%  '
$ws.Range('D41').Value = "'6.607"
$ws.Range('E41').Value = '  +3.48%  '
$ws.Range('D42').Value = "'0.9078"
$ws.Range('E42').Value = '  +0.31%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').Value = "'2.010.67"
$ws.Range('E44').Value = '  -18.29%  '
$ws.Range('D45').Value = "'101.52"
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = "'66.34"
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('D47').Value = "'1.755"
$ws.Range('E47').Value = '  +4.84%  '
$ws.Range('E48').Value = '  -1.15%  '
$ws.Range('E49').Value = '  +3.09%  '
$ws.Range('D50').Value = "'9.014"
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').Value = "'0.3954"
